# Sanity adjustment for "Semilla 9" sheet: update seed/test data values
# (phone numbers / port ids) and normalize formatting on row 14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semilla 9")

# --- Row 9 ---
$ws.Range("C9").Value = "3045981670"
$ws.Range("D9").Value = "732111324707276"

# --- Row 10 ---
$ws.Range("C10").Value = "3043209863"
$ws.Range("D10").Value = "732111324707277"

# --- Row 11 ---
$ws.Range("B11").Value = "309991475"
$ws.Range("C11").Value = "3045981670"
$ws.Range("D11").Value = "732111324707276"

# --- Row 12 ---
$ws.Range("B12").Value = "270670616"
$ws.Range("D12").Value = "732111193280551"

# --- Row 13 ---
$ws.Range("B13").Value = "163908584"
$ws.Range("D13").Value = "732111193280544"

# --- Row 14 ---
# Normalize the left-aligned style on B14/C14 to match the plain style used
# by the rest of the column (copy formatting from row 13, then restore values)
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122)

$ws.Range("B14").Value = "697979125"
$ws.Range("D14").Value = "732111193280535"

$ws.Range("C12").Value = "3046010569"
$ws.Range("C13").Value = "3046010523"
$ws.Range("C14").Value = "3046008593"
